$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Merging data" (row 40) is no longer listed in the Swedish sitemap (E column)
$ws.Range("E40").ClearContents()

# "Collapsing data" (row 41) lost its Swedish page/meta/sitemap tracking entirely
$ws.Range("C41:E41").ClearContents()

# "Glossary" heading (A13) is now bolded
$ws.Range("A13").Font.Bold = $true
